$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $cellRef, $val) {
    $r = $sheet.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "39.485.63"
$ws.Range("E2").Value = "  +1.87%  "
Set-TextValue $ws "D3" "2.158.02"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws "D5" "227.52"
$ws.Range("E5").Value = "  +0.04%  "
Set-TextValue $ws "D6" "0.626"
$ws.Range("E6").Value = "  +1.67%  "
Set-TextValue $ws "D7" "63.22"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue $ws "D9" "0.393"
$ws.Range("E9").Value = "  +0.89%  "
Set-TextValue $ws "D10" "0.0852"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +0.28%  "
Set-TextValue $ws "D12" "15.98"
$ws.Range("E12").Value = "  +1.68%  "
Set-TextValue $ws "D13" "2.479.87"
$ws.Range("E13").Value = "  +2.65%  "
Set-TextValue $ws "D14" "21.95"
$ws.Range("E14").Value = "  -0.29%  "
Set-TextValue $ws "D15" "0.810"
$ws.Range("E15").Value = "  -0.18%  "
Set-TextValue $ws "D16" "5.50"
$ws.Range("E16").Value = "  -0.59%  "
Set-TextValue $ws "D17" "2.155.87"
$ws.Range("E17").Value = "  +2.79%  "
Set-TextValue $ws "D18" "39.518.40"
$ws.Range("E18").Value = "  +1.89%  "
Set-TextValue $ws "D19" "72.19"
Set-TextValue $ws "D20" "6.15"
$ws.Range("E20").Value = "  +0.23%  "
Set-TextValue $ws "D21" "0.0₃0848"
$ws.Range("E21").Value = "  +0.63%  "
Set-TextValue $ws "D22" "229.24"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +1.23%  "
Set-TextValue $ws "D25" "2.26"
$ws.Range("E25").Value = "  -3.58%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D26" "9.66"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D27" "171.99"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +0.90%  "
Set-TextValue $ws "D29" "19.79"
$ws.Range("E29").Value = "  +2.37%  "
Set-TextValue $ws "D30" "1.41"
$ws.Range("E30").Value = "  -0.68%  "
Set-TextValue $ws "D31" "2.68"
$ws.Range("E31").Value = "  +5.35%  "
$ws.Range("E32").Value = "  +1.85%  "
Set-TextValue $ws "D33" "4.62"
$ws.Range("E33").Value = "  +1.51%  "
Set-TextValue $ws "D34" "4.68"
$ws.Range("E34").Value = "  -1.34%  "
Set-TextValue $ws "D35" "6.96"
$ws.Range("E35").Value = "  -1.72%  "
Set-TextValue $ws "D36" "0.0622"
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("E37").Value = "  +1.13%  "
Set-TextValue $ws "D38" "3.60"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("B40").Value = "FTXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D40" "4.71"
$ws.Range("E40").Value = "  +14.03%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D41" "102.26"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D42" "0.0227"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D43" "17.68"
$ws.Range("E43").Value = "  -2.26%  "
Set-TextValue $ws "D44" "1.518.09"
$ws.Range("E44").Value = "  -0.61%  "
Set-TextValue $ws "D45" "1.20"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E48").Value = "  -0.10%  "
Set-TextValue $ws "D49" "7.70"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D50" "2.99"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws "D51" "2.363.39"
$ws.Range("E51").Value = "  +2.62%  "
